# Apply the "Features_Placed" directory restructuring to the image map sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-14 hold the "placed feature" images for creature 1 (body, arms, legs, eyes 1-2).
# Their image paths move from Data/Images/... to Data/Images/Features_Placed/...
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    $updated = $current -replace '^Data/Images/', 'Data/Images/Features_Placed/'
    $cell.Value = $updated
}

# Widen column B to fit the longer paths and drop the old "best fit" auto-sizing.
$ws.Columns.Item(2).ColumnWidth = 40.17

# Update the selected/active cell to B14, matching the edited range.
$ws.Range("B14").Select()
